$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "2016-03-13 10:44:52"
$wsZh.Range("H3").Value = "2016-03-13 10:45:11"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "2016-03-13 10:44:56"
$wsDe.Range("H3").Value = "2016-03-13 10:45:17"
